$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates are Excel serial numbers, matching existing column A values)
$newRows = @(
    @{ Row = 230; A = 44304; B = 2; C = 27; D = 239.1708743024183 },
    @{ Row = 231; A = 44305; B = 0; C = 18; D = 159.4472495349455 },
    @{ Row = 232; A = 44306; B = 4; C = 22; D = 194.8799716538223 },
    @{ Row = 233; A = 44307; B = 0; C = 22; D = 194.8799716538223 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Column A: date serial, formatted like the existing date column (style copied from A229)
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Range("A229").Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) # xlPasteFormats

    # Columns B, C, D: plain numeric values
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
}
